$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "status"

$ws.Range("E2").Value = "Active"
$ws.Range("E3").Value = "Active"
$ws.Range("E4").Value = "Active"
$ws.Range("E7").Value = "Active"
$ws.Range("E8").Value = "Active"
$ws.Range("E11").Value = "Active"
$ws.Range("E13").Value = "Active"
$ws.Range("E14").Value = "Active"
$ws.Range("E19").Value = "Active"
$ws.Range("E20").Value = "Active"

$ws.Range("E6").Value = "Deactive"
$ws.Range("E9").Value = "Deactive"
$ws.Range("E10").Value = "Deactive"
$ws.Range("E16").Value = "Deactive"
$ws.Range("E17").Value = "Deactive"
$ws.Range("E18").Value = "Deactive"

$ws.Range("E5").Value = "Suspend"
$ws.Range("E12").Value = "Suspend"
$ws.Range("E15").Value = "Suspend"

$ws.Range("E6").Select()
